$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("servents")

# 1. C18: "Grandcaster" -> "Loregrandcaster"
$ws.Range("C18").Value = "Loregrandcaster"

# 2-11. Fill column F for rows 150-159 (previously empty inline strings)
$ws.Range("F150").Value = "幼體／迪亞馬特"
$ws.Range("F151").Value = "所多瑪之獸／德拉科"
$ws.Range("F152").Value = "果心居士"
$ws.Range("F153").Value = "怖軍"
$ws.Range("F154").Value = "杜爾迦"
$ws.Range("F155").Value = "美杜莎"
$ws.Range("F156").Value = "雨之魔女托內莉可"
$ws.Range("F157").Value = "阿爾托莉亞．Caster"
$ws.Range("F158").Value = "美露莘"
$ws.Range("F159").Value = "旺吉娜"

# 12-19. Fill column E for rows 164-171 (previously empty inline strings)
$ws.Range("E164").Value = "源賴光／醜御前"
$ws.Range("E165").Value = "安德洛墨達"
$ws.Range("E166").Value = "瑪麗·安託瓦內特〔Alter〕"
$ws.Range("E167").Value = "巖窟王　基督山"
$ws.Range("E168").Value = "蒼崎青子"
$ws.Range("E169").Value = "久遠寺有珠"
$ws.Range("E170").Value = "響＆千鍵"
$ws.Range("E171").Value = "埃列什基伽勒"

# 20. Add new rows 179-186 with data in columns A-D, and blank (but present) E/F cells
$newRows = @(
    @(433, 5, "Lancer", "ビショーネ"),
    @(435, 5, "Caster", "小野小町"),
    @(437, 5, "Pretender", "ダンテ・アリギエーリ"),
    @(438, 5, "Ruler", "メタトロン・ジャンヌ"),
    @(440, 5, "Berserker", "リリス"),
    @(441, 5, "Pretender", "テュフォン・エフェメロス"),
    @(442, 5, "Lancer", "インドラ"),
    @(444, 5, "'40", "Ｕ－オルガマリー")
)

$rowIdx = 179
foreach ($row in $newRows) {
    $ws.Cells.Item($rowIdx, 1).Value = $row[0]
    $ws.Cells.Item($rowIdx, 2).Value = $row[1]
    $ws.Cells.Item($rowIdx, 3).Value = $row[2]
    $ws.Cells.Item($rowIdx, 4).Value = $row[3]
    # Force creation of empty (but present) E and F cells, matching the
    # author's export which always emits all six columns per row.
    $ws.Cells.Item($rowIdx, 5).Style = "Normal"
    $ws.Cells.Item($rowIdx, 6).Style = "Normal"
    # Undo any quote-prefix / number-format styling the text assignments
    # above may have implicitly triggered, so the row keeps default style.
    $ws.Cells.Item($rowIdx, 3).Style = "Normal"
    $rowIdx++
}
